# =====================================================================
# Restructure PlayerPerformance_4677.xlsx:
#  1. Insert a new 'Player Info' sheet before 'ODI Batting'
#  2. On 'ODI Batting': MATCH_CARD_LINK -> MATCH_CODE (URL -> numeric
#     match code text), drop the stray empty INNING_NUMBER cells
#  3. On 'ODI Bowling': MATCH_CARD_LINK -> MATCH_CODE (URL -> numeric
#     match code text)
#  4. Append a new 'ODI Batting Extra' sheet after 'ODI Bowling'
# =====================================================================

$wb = $excel.ActiveWorkbook

function Set-HeaderStyle($rng) {
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}

$battingWs = $wb.Worksheets.Item('ODI Batting')
$bowlingWs = $wb.Worksheets.Item('ODI Bowling')

# --- 1) 'Player Info' sheet, inserted before 'ODI Batting' ---
$playerInfoWs = $wb.Worksheets.Add($battingWs)
$playerInfoWs.Name = 'Player Info'

$piHeaders = @('ID', 'NAME', 'BATTING_HAND', 'BOWL_STYLE')
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $cell = $playerInfoWs.Cells.Item(1, $c)
    $cell.Value = $piHeaders[$c - 1]
    Set-HeaderStyle($cell)
}

$playerInfoWs.Range('A2').NumberFormat = '@'
$playerInfoWs.Range('A2').Value = '4677'
$playerInfoWs.Range('B2').Value = 'Mohammad Saifuddin'
$playerInfoWs.Range('C2').Value = 'Left Handed'
$playerInfoWs.Range('D2').Value = 'Right Arm Medium Fast'

# --- 2) 'ODI Batting': MATCH_CARD_LINK header -> MATCH_CODE,
#        URL values -> bare numeric match-code text ---
# Re-resolve the reference: inserting 'Player Info' just above shifted
# sheet positions, and this COM binding tracks by position, not by the
# object originally fetched.
$battingWs = $wb.Worksheets.Item('ODI Batting')
$battingWs.Range('D1').Value = 'MATCH_CODE'

$battingWs.Range('D2:D30').NumberFormat = '@'
$battingWs.Range('D2').Value = '4080'
$battingWs.Range('D3').Value = '4116'
$battingWs.Range('D4').Value = '4124'
$battingWs.Range('D5').Value = '4214'
$battingWs.Range('D6').Value = '4217'
$battingWs.Range('D7').Value = '4218'
$battingWs.Range('D8').Value = '4230'
$battingWs.Range('D9').Value = '4250'
$battingWs.Range('D10').Value = '4251'
$battingWs.Range('D11').Value = '4252'
$battingWs.Range('D12').Value = '4286'
$battingWs.Range('D13').Value = '4295'
$battingWs.Range('D14').Value = '4296'
$battingWs.Range('D15').Value = '4307'
$battingWs.Range('D16').Value = '4311'
$battingWs.Range('D17').Value = '4314'
$battingWs.Range('D18').Value = '4325'
$battingWs.Range('D19').Value = '4335'
$battingWs.Range('D20').Value = '4345'
$battingWs.Range('D21').Value = '4349'
$battingWs.Range('D22').Value = '4416'
$battingWs.Range('D23').Value = '4420'
$battingWs.Range('D24').Value = '4447'
$battingWs.Range('D25').Value = '4453'
$battingWs.Range('D26').Value = '4463'
$battingWs.Range('D27').Value = '4464'
$battingWs.Range('D28').Value = '4477'
$battingWs.Range('D29').Value = '4479'
$battingWs.Range('D30').Value = '4481'

# Drop the stray empty INNING_NUMBER (column B) cells left over
# from the scrape -- these rows never had an inning number scraped.
$battingWs.Range('B6').Value = ''
$battingWs.Range('B7').Value = ''
$battingWs.Range('B8').Value = ''
$battingWs.Range('B12').Value = ''
$battingWs.Range('B13').Value = ''
$battingWs.Range('B14').Value = ''
$battingWs.Range('B15').Value = ''
$battingWs.Range('B18').Value = ''
$battingWs.Range('B23').Value = ''
$battingWs.Range('B30').Value = ''

# --- 3) 'ODI Bowling': MATCH_CARD_LINK header -> MATCH_CODE,
#        URL values -> bare numeric match-code text ---
# Re-resolve again for the same reason (position shifted by the earlier
# 'Player Info' insert).
$bowlingWs = $wb.Worksheets.Item('ODI Bowling')
$bowlingWs.Range('B1').Value = 'MATCH_CODE'

$bowlingWs.Range('B2:B29').NumberFormat = '@'
$bowlingWs.Range('B2').Value = '4080'
$bowlingWs.Range('B3').Value = '4116'
$bowlingWs.Range('B4').Value = '4124'
$bowlingWs.Range('B5').Value = '4214'
$bowlingWs.Range('B6').Value = '4217'
$bowlingWs.Range('B7').Value = '4218'
$bowlingWs.Range('B8').Value = '4230'
$bowlingWs.Range('B9').Value = '4250'
$bowlingWs.Range('B10').Value = '4251'
$bowlingWs.Range('B11').Value = '4252'
$bowlingWs.Range('B12').Value = '4286'
$bowlingWs.Range('B13').Value = '4295'
$bowlingWs.Range('B14').Value = '4296'
$bowlingWs.Range('B15').Value = '4307'
$bowlingWs.Range('B16').Value = '4311'
$bowlingWs.Range('B17').Value = '4314'
$bowlingWs.Range('B18').Value = '4325'
$bowlingWs.Range('B19').Value = '4335'
$bowlingWs.Range('B20').Value = '4345'
$bowlingWs.Range('B21').Value = '4349'
$bowlingWs.Range('B22').Value = '4416'
$bowlingWs.Range('B23').Value = '4420'
$bowlingWs.Range('B24').Value = '4447'
$bowlingWs.Range('B25').Value = '4453'
$bowlingWs.Range('B26').Value = '4463'
$bowlingWs.Range('B27').Value = '4477'
$bowlingWs.Range('B28').Value = '4479'
$bowlingWs.Range('B29').Value = '4481'

# --- 4) 'ODI Batting Extra' sheet, appended after 'ODI Bowling' ---
# Re-resolve the 'ODI Bowling' reference: inserting 'Player Info' earlier
# shifted sheet positions, and this COM binding tracks by position, not
# by the object originally fetched.
$bowlingWs = $wb.Worksheets.Item('ODI Bowling')
$extraWs = $wb.Worksheets.Add($null, $bowlingWs)
$extraWs.Name = 'ODI Batting Extra'

$exHeaders = @('MATCH_CODE', 'BATTING_POSITION', 'NUM_4', 'NUM_6', 'PERCENT_RUNS_OF_TOTAL', 'MAN_OF_MATCH')
for ($c = 1; $c -le $exHeaders.Length; $c++) {
    $cell = $extraWs.Cells.Item(1, $c)
    $cell.Value = $exHeaders[$c - 1]
    Set-HeaderStyle($cell)
}

# Columns A, C, D, E hold text (match codes / digit-strings / percent
# strings); column B (BATTING_POSITION) holds real numbers, column F
# (MAN_OF_MATCH) holds plain text.
$extraWs.Range('A2:A21').NumberFormat = '@'
$extraWs.Range('C2:E21').NumberFormat = '@'

$extraWs.Range('A2').Value = '4252'
$extraWs.Range('B2').Value = 7
$extraWs.Range('C2').Value = '4'
$extraWs.Range('D2').Value = '0'
$extraWs.Range('E2').Value = '18.18%'
$extraWs.Range('F2').Value = 'NO'

$extraWs.Range('A3').Value = '4286'
$extraWs.Range('B3').Value = 8
$extraWs.Range('C3').Value = ''
$extraWs.Range('D3').Value = ''
$extraWs.Range('E3').Value = ''
$extraWs.Range('F3').Value = 'NO'

$extraWs.Range('A4').Value = '4295'
$extraWs.Range('B4').Value = ''
$extraWs.Range('C4').Value = ''
$extraWs.Range('D4').Value = ''
$extraWs.Range('E4').Value = ''
$extraWs.Range('F4').Value = 'NO'

$extraWs.Range('A5').Value = '4296'
$extraWs.Range('B5').Value = 10
$extraWs.Range('C5').Value = ''
$extraWs.Range('D5').Value = ''
$extraWs.Range('E5').Value = ''
$extraWs.Range('F5').Value = 'NO'

$extraWs.Range('A6').Value = '4307'
$extraWs.Range('B6').Value = ''
$extraWs.Range('C6').Value = ''
$extraWs.Range('D6').Value = ''
$extraWs.Range('E6').Value = ''
$extraWs.Range('F6').Value = 'NO'

$extraWs.Range('A7').Value = '4311'
$extraWs.Range('B7').Value = ''
$extraWs.Range('C7').Value = ''
$extraWs.Range('D7').Value = ''
$extraWs.Range('E7').Value = ''
$extraWs.Range('F7').Value = 'NO'

$extraWs.Range('A8').Value = '4314'
$extraWs.Range('B8').Value = 8
$extraWs.Range('C8').Value = '0'
$extraWs.Range('D8').Value = '0'
$extraWs.Range('E8').Value = '1.79%'
$extraWs.Range('F8').Value = 'NO'

$extraWs.Range('A9').Value = '4325'
$extraWs.Range('B9').Value = ''
$extraWs.Range('C9').Value = ''
$extraWs.Range('D9').Value = ''
$extraWs.Range('E9').Value = ''
$extraWs.Range('F9').Value = 'NO'

$extraWs.Range('A10').Value = '4335'
$extraWs.Range('B10').Value = 8
$extraWs.Range('C10').Value = '0'
$extraWs.Range('D10').Value = '0'
$extraWs.Range('E10').Value = '0.76%'
$extraWs.Range('F10').Value = 'NO'

$extraWs.Range('A11').Value = '4345'
$extraWs.Range('B11').Value = 8
$extraWs.Range('C11').Value = '9'
$extraWs.Range('D11').Value = '0'
$extraWs.Range('E11').Value = '17.83%'
$extraWs.Range('F11').Value = 'NO'

$extraWs.Range('A12').Value = '4349'
$extraWs.Range('B12').Value = 8
$extraWs.Range('C12').Value = '0'
$extraWs.Range('D12').Value = '0'
$extraWs.Range('E12').Value = ''
$extraWs.Range('F12').Value = 'NO'

$extraWs.Range('A13').Value = '4416'
$extraWs.Range('B13').Value = 7
$extraWs.Range('C13').Value = '0'
$extraWs.Range('D13').Value = '3'
$extraWs.Range('E13').Value = '8.72%'
$extraWs.Range('F13').Value = 'NO'

$extraWs.Range('A14').Value = '4420'
$extraWs.Range('B14').Value = 10
$extraWs.Range('C14').Value = ''
$extraWs.Range('D14').Value = ''
$extraWs.Range('E14').Value = ''
$extraWs.Range('F14').Value = 'NO'

$extraWs.Range('A15').Value = '4447'
$extraWs.Range('B15').Value = 8
$extraWs.Range('C15').Value = '1'
$extraWs.Range('D15').Value = '0'
$extraWs.Range('E15').Value = '1.68%'
$extraWs.Range('F15').Value = 'NO'

$extraWs.Range('A16').Value = '4453'
$extraWs.Range('B16').Value = 8
$extraWs.Range('C16').Value = '0'
$extraWs.Range('D16').Value = '0'
$extraWs.Range('E16').Value = '2.58%'
$extraWs.Range('F16').Value = 'NO'

$extraWs.Range('A17').Value = '4463'
$extraWs.Range('B17').Value = 8
$extraWs.Range('C17').Value = '2'
$extraWs.Range('D17').Value = '0'
$extraWs.Range('E17').Value = '5.06%'
$extraWs.Range('F17').Value = 'NO'

$extraWs.Range('A18').Value = '4464'
$extraWs.Range('B18').Value = ''
$extraWs.Range('C18').Value = ''
$extraWs.Range('D18').Value = ''
$extraWs.Range('E18').Value = ''
$extraWs.Range('F18').Value = 'NO'

$extraWs.Range('A19').Value = '4477'
$extraWs.Range('B19').Value = 9
$extraWs.Range('C19').Value = '1'
$extraWs.Range('D19').Value = '0'
$extraWs.Range('E19').Value = '2.90%'
$extraWs.Range('F19').Value = 'NO'

$extraWs.Range('A20').Value = '4479'
$extraWs.Range('B20').Value = 9
$extraWs.Range('C20').Value = '1'
$extraWs.Range('D20').Value = '0'
$extraWs.Range('E20').Value = '11.57%'
$extraWs.Range('F20').Value = 'NO'

$extraWs.Range('A21').Value = '4481'
$extraWs.Range('B21').Value = 9
$extraWs.Range('C21').Value = ''
$extraWs.Range('D21').Value = ''
$extraWs.Range('E21').Value = ''
$extraWs.Range('F21').Value = 'NO'

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $s = $wb.Worksheets.Item($i)
    Write-Host "Sheet $i : $($s.Name)"
}

